$wb = $excel.ActiveWorkbook
$win = $excel.ActiveWindow
Write-Host $win.WindowState
try { Write-Host $win.Left } catch {}
try { Write-Host $win.Top } catch {}
